$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "[Rogério-Processos de Usinagem 2, -, -, -]"
$ws.Range("C2").Value = "[Victor Lima-CAM, -, Emerson-Robótica, Rogério-Processos de Usinagem 2]"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "[-, -, João Paulo-Lab. de eletroeletrônica]"

# Row 3
$ws.Range("B3").Value = "[Rogério-Processos de Usinagem 2, -, -, -]"
$ws.Range("C3").Value = "[Ludoff-Eletrohidráulica, Pedro Bispo-Automação Industrial, Nilton Maia-Elementos de máquinas, Leonardo-Manut. Mecânica]"
$ws.Range("D3").Value = "Cleidson-Máquinas Elétri"
$ws.Range("E3").Value = "[Euclides-Tecnologia da soldagem, Humberto-Eletropneumática, Eudes-Microcontroladores, Carlos Eduardo-Processos de Usinagem 1]"
$ws.Range("F3").Value = "-"

# Row 4
$ws.Range("B4").Value = "[-, Pedro Bispo-Lab. Máquinas Elétricas, Pedro Bispo-Lab. Máquinas Elétricas]"
$ws.Range("C4").Value = "[Ludoff-Eletrohidráulica, Pedro Bispo-Automação Industrial, Nilton Maia-Elementos de máquinas, Leonardo-Manut. Mecânica]"
$ws.Range("D4").Value = "Cleidson-Máquinas Elétri"
$ws.Range("E4").Value = "[Euclides-Tecnologia da soldagem, Humberto-Eletropneumática, Eudes-Microcontroladores, Carlos Eduardo-Processos de Usinagem 1]"
$ws.Range("F4").Value = "-"

# Row 6
$ws.Range("B6").Value = "[Emerson-Robótica, -, -, -]"
$ws.Range("C6").Value = "[Ludoff-Eletrohidráulica, Pedro Bispo-Automação Industrial, Nilton Maia-Elementos de máquinas, Leonardo-Manut. Mecânica]"
$ws.Range("D6").Value = "André Guimarães-Máquinas Térmicas e de Fl"
$ws.Range("E6").Value = "[Euclides-Tecnologia da soldagem, Humberto-Eletropneumática, Eudes-Microcontroladores, Carlos Eduardo-Processos de Usinagem 1]"
$ws.Range("F6").Value = "[-, -, Victor Lima-CAM, -]"

# Row 7
$ws.Range("B7").Value = "[Emerson-Robótica, -, -, -]"
$ws.Range("C7").Value = "[Ludoff-Eletrohidráulica, Pedro Bispo-Automação Industrial, Nilton Maia-Elementos de máquinas, Leonardo-Manut. Mecânica]"
$ws.Range("D7").Value = "André Guimarães-Máquinas Térmicas e de Fl"
$ws.Range("E7").Value = "[Euclides-Tecnologia da soldagem, Humberto-Eletropneumática, Eudes-Microcontroladores, Carlos Eduardo-Processos de Usinagem 1]"
$ws.Range("F7").Value = "-"

# Row 8
$ws.Range("B8").Value = "[Victor Lima-CAM, -, -, -]"
$ws.Range("C8").Value = "[Rogério-Processos de Usinagem 2, -, Victor Lima-CAM, Emerson-Robótica]"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "[-, -, João Paulo-Lab. de eletroeletrônica]"
